$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '246.84'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '0.58%'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '29.56'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '8.71%'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '5.167'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '1.18%'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.05728'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '0.42%'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '6.579'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '1.08%'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.098'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '2.97%'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.8569'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '4.61%'
$ws.Range('B9').Value = 'FTXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.8671'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '0.95%'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.1367'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '2.68%'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07084'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '2.14%'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.02931'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '3.06%'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.09379'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '-0.16%'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.001523'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '-0.37%'
$ws.Range('B15').Value = 'CoinExToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.04135'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '2.44%'
$ws.Range('B16').Value = 'One'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0006007'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '-94.06%'
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.006122'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '-1.44%'
$ws.Range('B18').Value = 'UpBots'
$ws.Range('C18').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.007489'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '5,071.37%'
$ws.Range('B19').Value = 'LEO'
$ws.Range('C19').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.489'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '-0.62%'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '2.279'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '2.23%'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.3182'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '0.54%'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.03395'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '5.27%'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '2.29%'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.476'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '-2.78%'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1379'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '0.44%'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.005008'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '11.96%'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.001225'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '0.49%'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '22.27%'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.03753'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '0.71%'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.005767'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '66.95%'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1072'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '1.17%'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.002427'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '-0.46%'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.008491'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '-7.91%'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00005246'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '1.69%'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '0.00%'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.06467'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '-35.94%'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.002533'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '1.13%'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '0.00%'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '0.00%'
